$wb = $excel.ActiveWorkbook

# Sheet "Admin": update cell D2 with the new value
$wsAdmin = $wb.Worksheets.Item("Admin")
$wsAdmin.Range("D2").Value = "8448782A"

# Sheet "Jira": update rows 3 and 4
$wsJira = $wb.Worksheets.Item("Jira")
$wsJira.Range("A3").Value = "Leave_Accept"
$wsJira.Range("A4").Value = "Recruitment_Rejected1"

# B4 becomes an empty text cell (same as B3, which already holds an empty string)
$wsJira.Range("B4").Value = "'"
$wsJira.Range("B4").Style = $wsJira.Range("B3").Style
